$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value2 = 10.12586833333333
$ws.Range("H2").Value2 = 30.377605
$ws.Range("I2").Value2 = 0.9311967029481902
$ws.Range("J2").Value2 = 0.9311967029481902
$ws.Range("M2").Value2 = 48.42420966666666
$ws.Range("N2").Value2 = 145.272629
$ws.Range("O2").Value2 = 0.6311762527593259
$ws.Range("P2").Value2 = 0.6311762527593258
$ws.Range("Q2").Value2 = 490.3371712303939
$ws.Range("R2").Value2 = 4413.034541073545
$ws.Range("S2").Value2 = 0.5877492455486778
$ws.Range("T2").Value2 = 0.5877492455486777
$ws.Range("G3").Value2 = 10.12586833333333
$ws.Range("H3").Value2 = 30.377605
$ws.Range("I3").Value2 = 0.9311967029481902
$ws.Range("J3").Value2 = 0.9311967029481902
$ws.Range("M3").Value2 = 6.849914666666667
$ws.Range("O3").Value2 = 0.08928392431779728
$ws.Range("P3").Value2 = 0.08928392431779726
$ws.Range("Q3").Value2 = 69.36133400923556
$ws.Range("R3").Value2 = 624.25200608312
$ws.Range("S3").Value2 = 0.08314089595100857
$ws.Range("T3").Value2 = 0.08314089595100856
$ws.Range("G4").Value2 = 10.12586833333333
$ws.Range("H4").Value2 = 30.377605
$ws.Range("I4").Value2 = 0.9311967029481902
$ws.Range("J4").Value2 = 0.9311967029481902
$ws.Range("N4").Value2 = 64.33937399999999
$ws.Range("O4").Value2 = 0.2795398229228769
$ws.Range("P4").Value2 = 0.2795398229228769
$ws.Range("Q4").Value2 = 217.1640099243633
$ws.Range("R4").Value2 = 1954.47608931927
$ws.Range("S4").Value2 = 0.2603065614485039
$ws.Range("T4").Value2 = 0.2603065614485039
$ws.Range("I5").Value2 = 0.009287810103293732
$ws.Range("J5").Value2 = 0.009287810103293733
$ws.Range("M5").Value2 = 48.42420966666666
$ws.Range("N5").Value2 = 145.272629
$ws.Range("O5").Value2 = 0.6311762527593259
$ws.Range("P5").Value2 = 0.6311762527593258
$ws.Range("Q5").Value2 = 4.890651479494665
$ws.Range("R5").Value2 = 44.01586331545199
$ws.Range("S5").Value2 = 0.005862245177337145
$ws.Range("T5").Value2 = 0.005862245177337145
$ws.Range("I6").Value2 = 0.009287810103293732
$ws.Range("J6").Value2 = 0.009287810103293733
$ws.Range("M6").Value2 = 6.849914666666667
$ws.Range("O6").Value2 = 0.08928392431779728
$ws.Range("P6").Value2 = 0.08928392431779726
$ws.Range("Q6").Value2 = 0.6918139816746666
$ws.Range("S6").Value2 = 0.0008292521343405505
$ws.Range("T6").Value2 = 0.0008292521343405505
$ws.Range("I7").Value2 = 0.009287810103293732
$ws.Range("J7").Value2 = 0.009287810103293733
$ws.Range("N7").Value2 = 64.33937399999999
$ws.Range("O7").Value2 = 0.2795398229228769
$ws.Range("P7").Value2 = 0.2795398229228769
$ws.Range("Q7").Value2 = 2.166006472167999
$ws.Range("S7").Value2 = 0.002596312791616037
$ws.Range("T7").Value2 = 0.002596312791616037
$ws.Range("I8").Value2 = 0.05951548694851595
$ws.Range("J8").Value2 = 0.05951548694851596
$ws.Range("M8").Value2 = 48.42420966666666
$ws.Range("N8").Value2 = 145.272629
$ws.Range("O8").Value2 = 0.6311762527593259
$ws.Range("P8").Value2 = 0.6311762527593258
$ws.Range("Q8").Value2 = 31.33887332541211
$ws.Range("R8").Value2 = 282.049859928709
$ws.Range("S8").Value2 = 0.03756476203331086
$ws.Range("T8").Value2 = 0.03756476203331086
$ws.Range("I9").Value2 = 0.05951548694851595
$ws.Range("J9").Value2 = 0.05951548694851596
$ws.Range("M9").Value2 = 6.849914666666667
$ws.Range("O9").Value2 = 0.08928392431779728
$ws.Range("P9").Value2 = 0.08928392431779726
$ws.Range("Q9").Value2 = 4.433084391180445
$ws.Range("R9").Value2 = 39.897759520624
$ws.Range("S9").Value2 = 0.00531377623244815
$ws.Range("T9").Value2 = 0.005313776232448149
$ws.Range("I10").Value2 = 0.05951548694851595
$ws.Range("J10").Value2 = 0.05951548694851596
$ws.Range("N10").Value2 = 64.33937399999999
$ws.Range("O10").Value2 = 0.2795398229228769
$ws.Range("P10").Value2 = 0.2795398229228769
$ws.Range("Q10").Value2 = 13.87958286087266
$ws.Range("S10").Value2 = 0.01663694868275694
$ws.Range("T10").Value2 = 0.01663694868275694
